$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shipment rows to append below the existing data (rows 6-13).
$rows = @(
    @("240037217222", "JESUSU ALBERTO HINCAPIE PJ", "3173734168", "$ 267.000,00"),
    @("240036364549", "Pedro Antonio Gomez", "3043925139", "$ 37.000,00"),
    @("240037128792", "EDGAR BUITRAGO", "3246465852", "$ 80.000,00"),
    @("240037108787", "Martha Apolonia Galvis Portillo", "3157200511", "$ 129.900,00"),
    @("700170049543", "RINA ISABEL BERMUDEZ GUERRA", "3226614162", "$ 45.000,00"),
    @("700170495277", "RINA BERMUDEZ GUERRA", "3226614162", "$ 45.000,00"),
    @("240037036846", "Héctor Fabio Bastidas", "3152974198", "$ 94.900,00"),
    @("240037080712", "Héctor Fabio Bastidas", "3152974198", "$ 100.000,00")
)

$r = 6
foreach ($row in $rows) {
    $ws.Range("A$r").NumberFormat = "@"
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").NumberFormat = "@"
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $r++
}
